$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 6.48142807727062, 0, 14.40014219143469)
    3 = @(0.3464964993005633, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 3.21235251628842)
    4 = @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 0, 5.488907176552729)
    5 = @(0.02258322285507441, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 2.888439239842931)
    6 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 6.048734245549538)
    7 = @(0.7287194209349384, 0.3375848360084654, 0.1529057820181812, 0.4998867070740569, 1, 1.719096746035642)
    8 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 6.048734245549538)
    9 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 6.048734245549538)
}

$cols = @("B", "C", "D", "E", "F", "G")

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $row).Value = $vals[$i]
    }
}
